$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A51").Value = "G4"
$ws.Range("B51").Value = "Read Book"
$ws.Range("C51").Value = $ws.Range("C50").Value2()
$ws.Range("C51").NumberFormat = $ws.Range("C50").NumberFormat()
$ws.Range("D51").Value = 1
$ws.Range("E51").Value = 0
$ws.Range("F51").Value = 0
